# Apply weekly cryptocurrency price/volume refresh to the "cryptos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never let Excel auto-coerce
# dotted price strings like "30.482.61" or "1.002" into numbers).
# We stage the text in a scratch cell via a text-returning formula,
# copy it, and paste-special (values only) into the destination so the
# destination cell keeps its original (unstyled) formatting.
function Set-TextValue {
    param($cellAddr, [string]$text)
    $scratch = $ws.Range("Z100")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($cellAddr).PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
}

Set-TextValue "D2" "30.482.61"
Set-TextValue "E2" "  -1.53%  "

Set-TextValue "D3" "2.097.60"
Set-TextValue "E3" "  -0.88%  "

Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.09%  "

Set-TextValue "D5" "330.35"
Set-TextValue "E5" "  -1.28%  "

Set-TextValue "E6" "  -0.09%  "

Set-TextValue "D7" "0.5223"
Set-TextValue "E7" "  -2.12%  "

Set-TextValue "D8" "0.4432"
Set-TextValue "E8" "  +0.92%  "

Set-TextValue "D9" "53.61"
Set-TextValue "E9" "  +15.97%  "

Set-TextValue "D10" "0.08934"
Set-TextValue "E10" "  -1.38%  "

Set-TextValue "E11" "  -2.37%  "

Set-TextValue "E12" "  -3.26%  "

Set-TextValue "D13" "2.091.03"
Set-TextValue "E13" "  -1.28%  "

Set-TextValue "D14" "6.692"
Set-TextValue "E14" "  -1.50%  "

Set-TextValue "D15" "7.709"
Set-TextValue "E15" "  -1.30%  "

Set-TextValue "E16" "  -1.32%  "

Set-TextValue "D17" "1.002"
Set-TextValue "E17" "  -0.22%  "

Set-TextValue "E18" "  -1.34%  "

Set-TextValue "D19" "0.06613"
Set-TextValue "E19" "  -0.97%  "

Set-TextValue "D20" "19.14"
Set-TextValue "E20" "  -0.35%  "

Set-TextValue "D22" "6.279"
Set-TextValue "E22" "  -1.60%  "

Set-TextValue "D23" "30.519.92"
Set-TextValue "E23" "  -1.65%  "

Set-TextValue "D24" "12.28"

Set-TextValue "D25" "2.319"

Set-TextValue "D26" "2.335.35"
Set-TextValue "E26" "  -1.27%  "

Set-TextValue "D27" "22.27"

Set-TextValue "D28" "2.571"

Set-TextValue "D29" "163.68"
Set-TextValue "E29" "  +0.06%  "

Set-TextValue "D30" "132.03"
Set-TextValue "E30" "  -1.40%  "

Set-TextValue "D31" "1.191"
Set-TextValue "E31" "  +1.31%  "

Set-TextValue "D32" "0.1071"
Set-TextValue "E32" "  -0.42%  "

Set-TextValue "D33" "1.661"
Set-TextValue "E33" "  +8.36%  "

Set-TextValue "D34" "6.164"
Set-TextValue "E34" "  -1.35%  "

Set-TextValue "D35" "3.897"
Set-TextValue "E35" "  -2.85%  "

Set-TextValue "D36" "10.20"
Set-TextValue "E36" "  +6.85%  "

Set-TextValue "D37" "0.02560"
Set-TextValue "E37" "  -2.32%  "

Set-TextValue "D38" "0.06801"
Set-TextValue "E38" "  +0.72%  "

Set-TextValue "B39" "Aptos"
Set-TextValue "C39" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D39" "12.74"
Set-TextValue "E39" "  -2.34%  "

Set-TextValue "B40" "InternetComputer(DFINITY)"
Set-TextValue "C40" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D40" "5.469"
Set-TextValue "E40" "  -1.78%  "

Set-TextValue "D41" "0.2261"
Set-TextValue "E41" "  -0.96%  "

Set-TextValue "D42" "0.6908"
Set-TextValue "E42" "  +0.21%  "

Set-TextValue "D43" "1.253"
Set-TextValue "E43" "  -0.18%  "

Set-TextValue "D44" "1.000"
Set-TextValue "E44" "  -0.09%  "

Set-TextValue "D45" "13.95"
Set-TextValue "E45" "  -1.21%  "

Set-TextValue "D46" "0.6342"
Set-TextValue "E46" "  -2.37%  "

Set-TextValue "E47" "  -2.12%  "

Set-TextValue "E48" "  -1.52%  "

Set-TextValue "D49" "1.245"
Set-TextValue "E49" "  +6.66%  "

Set-TextValue "E50" "  -2.78%  "

Set-TextValue "D51" "81.87"
Set-TextValue "E51" "  -1.60%  "
